# KAHI Session 2 Update
# Fill in the "Post Treatment" style self-rating answers (column B, rows 2-15)
# on the survey sheet, then leave the selection on B16 (the next empty cell)
# to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "A little worse "
$ws.Range("B3").Value  = "Somewhat worse "
$ws.Range("B4").Value  = "Somewhat worse "
$ws.Range("B5").Value  = "Somewhat worse "
$ws.Range("B6").Value  = "Somewhat worse "
$ws.Range("B7").Value  = "A lot worse "
$ws.Range("B8").Value  = "Somewhat worse "
$ws.Range("B9").Value  = "Somewhat worse "
$ws.Range("B10").Value = "A little worse "
$ws.Range("B11").Value = "A lot worse "
$ws.Range("B12").Value = "A little worse "
$ws.Range("B13").Value = "A lot worse "
$ws.Range("B14").Value = "Somewhat worse "
$ws.Range("B15").Value = "A little worse "

# Move the active selection to B16, matching the workbook's saved cursor
# position after data entry.
[void]$ws.Range("B16").Select()
